$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-8
# from serial 45233 (2023-11-03) to serial 45243 (2023-11-13)
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45243
}
